# Auto-generated PowerShell COM-interop script to apply the HPL Internship Daily Tasks edit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Text blocks for new task descriptions (shared strings) ----
$s62 = @"
1. Mailed test documents(Test Execution Report, Defect Report) to Team
2. Mailed MOF Entry Remarks to Khaled Sir
3. Got new task to test SAP Modules - Assigned by Khaled Sir
4. Started exploring Jira - Zephyr
"@
$s63 = @"
1. Exploring Jira - Zephyr (Imported Manual Test Cases into Zephyr)
2. Worked on JS Tasks
"@
$s64 = @"
1. Worked on JS Tasks
2. Worked on Zephyr
3. Khaled Sir will provide update on SAP 
"@
$s65 = @"
1. Worked on JS Tasks
2. Tested SAP - Product Allocation Module
3. New Assignment - Shahiduzzaman Sir - HCL Website - Deadline - 15-10-25
"@
$s66 = @"
1. Worked on HCL Website Assignment
"@
$s67 = @"
1. HCL Website Assignment
2. Completed SAP - Product Allocation Module Test
3. Completed Phase - 1 of Roadmap.(Completing Watching JS Videos - 90+ Completed) (Not Submitted to Anik Vaiya Yet.)
"@
$s68 = @"
1. Working on HCL Website
2. Starting Phase - 2 of the Roadmap
"@
$s69 = @"
1. HCL Website Assignment
2. Started Phase - 2 of the Roadmap
"@
$s70 = @"
1. HCL Website Update - 
(About Us , Mission & Vission , Meet Our Leaders)
2. Talked with Nayeem Vai - First make demo, then resources will be given. 
Take References from HPL , bdhfl
"@
$s73 = @"
1. HCL Website Update - 
(About Us , Mission & Vission , Meet Our Leaders)
"@
$LEAVE = @"
LEAVE
"@
$VACATION = @"
VACATION
"@

# ---- Step 1: copy cell formatting (style indices) from existing rows with matching layout ----
# Row 105 (A:date-style, B:day-style, C:wrap-text style) is a good source for the common (2,3,8) combo
$ws.Range("A105:C105").Copy() | Out-Null
$ws.Range("A106:C109").PasteSpecial(-4122) | Out-Null
$ws.Range("A111:C112").PasteSpecial(-4122) | Out-Null
$ws.Range("A115:C117").PasteSpecial(-4122) | Out-Null

# Row 106-109 before edit already used style 4 on C (cstyle match for row 110); keep by pasting A/B only then restyle C for 110
$ws.Range("A105:B105").Copy() | Out-Null
$ws.Range("A110:B110").PasteSpecial(-4122) | Out-Null
$ws.Range("A113:B114").PasteSpecial(-4122) | Out-Null
$ws.Range("A118:B121").PasteSpecial(-4122) | Out-Null
$ws.Range("A122:B126").PasteSpecial(-4122) | Out-Null

# C column style 6 (bold red, centered) currently only exists at B18 - copy format from there
$ws.Range("B18").Copy() | Out-Null
$ws.Range("C113:C114").PasteSpecial(-4122) | Out-Null
$ws.Range("C118:C121").PasteSpecial(-4122) | Out-Null

# C column style 7 (rows 122-126) already default there; ensure no stray formatting remains (no-op, already correct)

# ---- Step 2: set values ----
$ws.Range("A106").Value = 45915
$ws.Range("B106").Value = "Monday"
$ws.Range("C106").Value = $s62
$ws.Range("A107").Value = 45916
$ws.Range("B107").Value = "Tuesday"
$ws.Range("C107").Value = $s63
$ws.Range("A108").Value = 45917
$ws.Range("B108").Value = "Wednesday"
$ws.Range("C108").Value = $s64
$ws.Range("A109").Value = 45918
$ws.Range("B109").Value = "Thursday"
$ws.Range("C109").Value = $s65
$ws.Range("A110").Value = 45921
$ws.Range("B110").Value = "Sunday"
$ws.Range("C110").Value = $s66
$ws.Range("A111").Value = 45922
$ws.Range("B111").Value = "Monday"
$ws.Range("C111").Value = $s67
$ws.Range("A112").Value = 45923
$ws.Range("B112").Value = "Tuesday"
$ws.Range("C112").Value = $s69
$ws.Range("A113").Value = 45924
$ws.Range("B113").Value = "Wednesday"
$ws.Range("C113").Value = $LEAVE
$ws.Range("A114").Value = 45925
$ws.Range("B114").Value = "Thursday"
$ws.Range("C114").Value = $LEAVE
$ws.Range("A115").Value = 45928
$ws.Range("B115").Value = "Sunday"
$ws.Range("C115").Value = $s68
$ws.Range("A116").Value = 45929
$ws.Range("B116").Value = "Monday"
$ws.Range("C116").Value = $s70
$ws.Range("A117").Value = 45930
$ws.Range("B117").Value = "Tuesday"
$ws.Range("C117").Value = $s73
$ws.Range("A118").Value = 45931
$ws.Range("B118").Value = "Wednesday"
$ws.Range("C118").Value = $VACATION
$ws.Range("A119").Value = 45932
$ws.Range("B119").Value = "Thursday"
$ws.Range("C119").Value = $VACATION
$ws.Range("A120").Value = 45933
$ws.Range("B120").Value = "Friday"
$ws.Range("C120").Value = $VACATION
$ws.Range("A121").Value = 45934
$ws.Range("B121").Value = "Saturday"
$ws.Range("C121").Value = $VACATION
$ws.Range("A122").Value = 45935
$ws.Range("B122").Value = "Sunday"
$ws.Range("A123").Value = 45936
$ws.Range("B123").Value = "Monday"
$ws.Range("A124").Value = 45937
$ws.Range("B124").Value = "Tuesday"
$ws.Range("A125").Value = 45938
$ws.Range("B125").Value = "Wednesday"
$ws.Range("A126").Value = 45939
$ws.Range("B126").Value = "Thursday"

# ---- Step 3: row heights ----
$ws.Range("A106").RowHeight = 211.5
$ws.Range("A107").RowHeight = 105.75
$ws.Range("A108").RowHeight = 114.75
$ws.Range("A109").RowHeight = 171
$ws.Range("A110").RowHeight = 36
$ws.Range("A111").RowHeight = 195
$ws.Range("A112").RowHeight = 70.5
$ws.Range("A113").RowHeight = 36
$ws.Range("A114").RowHeight = 36
$ws.Range("A115").RowHeight = 70.5
$ws.Range("A116").RowHeight = 176.25
$ws.Range("A117").RowHeight = 70.5
$ws.Range("A118").RowHeight = 36
$ws.Range("A119").RowHeight = 36
$ws.Range("A120").RowHeight = 36
$ws.Range("A121").RowHeight = 36
$ws.Range("A122").RowHeight = 36
$ws.Range("A123").RowHeight = 36
$ws.Range("A124").RowHeight = 36
$ws.Range("A125").RowHeight = 36
$ws.Range("A126").RowHeight = 36
$ws.Range("A127").RowHeight = 36
$ws.Range("A128").RowHeight = 36

# ---- Step 4: rows 127/128 column B style change (s=7 -> s=3, values stay empty) ----
$ws.Range("B106").Copy() | Out-Null
$ws.Range("B127:B128").PasteSpecial(-4122) | Out-Null

# ---- Step 5: sheet view changes ----
$ws.Range("C8:J8").Select() | Out-Null
